# Applies the "Updated cryptos list" diff to before.xlsx.
# Source cells are plain text (t="inlineStr") in the workbook; the diff only
# rewrites <t> text content for D (Price) / E (Volume(1h)) and, for the three
# rows that got re-ranked/re-labelled, B (Coin) and C (Link) as well.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values that parse as a plain number (single decimal point) would be
# auto-converted to a numeric cell by Excel; prefix those with an apostrophe
# so they stay text, exactly like the original "##.##" style price strings.

$ws.Range("D2").Value = "41.918.80"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").Value = "2.290.91"
$ws.Range("E3").Value = "  -3.61%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'316.90"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "'102.78"
$ws.Range("E6").Value = "  -5.39%  "
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  -3.53%  "
$ws.Range("E10").Value = "  -7.08%  "
$ws.Range("E11").Value = "  -2.89%  "
$ws.Range("E12").Value = "  -4.06%  "
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "'0.963"
$ws.Range("E14").Value = "  -5.16%  "
$ws.Range("D15").Value = "'15.20"
$ws.Range("E15").Value = "  -6.35%  "
$ws.Range("D16").Value = "2.633.97"
$ws.Range("E16").Value = "  -3.59%  "
$ws.Range("D17").Value = "2.277.62"
$ws.Range("E17").Value = "  -4.89%  "
$ws.Range("D18").Value = "42.041.90"
$ws.Range("E18").Value = "  -2.10%  "
$ws.Range("D19").Value = "'7.46"
$ws.Range("E19").Value = "  -3.30%  "
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").Value = "'3.65"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "'73.34"
$ws.Range("E22").Value = "  -4.01%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'282.35"
$ws.Range("E23").Value = "  +9.22%  "
$ws.Range("D24").Value = "'9.99"
$ws.Range("E24").Value = "  +5.71%  "
$ws.Range("E25").Value = "  -3.75%  "
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("E27").Value = "  -6.43%  "
$ws.Range("D28").Value = "'2.34"
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("D29").Value = "'22.95"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").Value = "'35.61"
$ws.Range("E30").Value = "  -4.25%  "
$ws.Range("D31").Value = "'163.46"
$ws.Range("E31").Value = "  -5.24%  "
$ws.Range("D32").Value = "'0.0875"
$ws.Range("E32").Value = "  -2.78%  "
$ws.Range("D33").Value = "'5.83"
$ws.Range("E33").Value = "  -4.44%  "
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("D36").Value = "'0.114"
$ws.Range("E36").Value = "  -5.18%  "
$ws.Range("D37").Value = "'4.52"
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("D38").Value = "'2.89"
$ws.Range("E38").Value = "  +6.68%  "
$ws.Range("E39").Value = "  -5.12%  "
$ws.Range("D40").Value = "'3.66"
$ws.Range("E40").Value = "  -7.35%  "
$ws.Range("D41").Value = "'99.82"
$ws.Range("E41").Value = "  +11.55%  "
$ws.Range("D42").Value = "'1.46"
$ws.Range("E42").Value = "  -4.83%  "
$ws.Range("D43").Value = "'69.56"
$ws.Range("E43").Value = "  -3.40%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  -7.58%  "
$ws.Range("D46").Value = "'114.29"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").Value = "'11.87"
$ws.Range("E47").Value = "  -3.97%  "
$ws.Range("D48").Value = "'76.89"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").Value = "'8.94"
$ws.Range("E49").Value = "  -3.42%  "
$ws.Range("E50").Value = "  -5.56%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.25"
$ws.Range("E51").Value = "  -4.34%  "

Write-Host "Applied 88 cell updates"
